$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are written as literal text (matching the source inlineStr cells).
# Purely-numeric-looking strings are given a leading apostrophe so Excel
# stores them as text instead of silently converting to a Number.

$ws.Range("D2").Value = '30.860.22'
$ws.Range("E2").Value = '  -1.41%  '
$ws.Range("D3").Value = '1.941.99'
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("E4").Value = '  -0.35%  '
$ws.Range("D5").Value = '''243.14'
$ws.Range("E5").Value = '  -1.06%  '
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").Value = '''0.4904'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '''0.2934'
$ws.Range("E8").Value = '  -1.55%  '
$ws.Range("D9").Value = '''0.06898'
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").Value = '''19.26'
$ws.Range("E10").Value = '  -0.48%  '
$ws.Range("D11").Value = '''105.58'
$ws.Range("E11").Value = '  -2.15%  '
$ws.Range("D12").Value = '1.953.45'
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("D13").Value = '''0.07763'
$ws.Range("E13").Value = '  +0.02%  '
$ws.Range("D14").Value = '''5.362'
$ws.Range("E14").Value = '  -1.85%  '
$ws.Range("D15").Value = '''0.7038'
$ws.Range("E15").Value = '  -1.39%  '
$ws.Range("D16").Value = '''276.73'
$ws.Range("E16").Value = '  -3.48%  '
$ws.Range("D17").Value = '30.882.03'
$ws.Range("E17").Value = '  -1.15%  '
$ws.Range("D18").Value = '''0.000007736'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").Value = '''13.10'
$ws.Range("E19").Value = '  -0.97%  '
$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").Value = '2.207.21'
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '''1.001'
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("D22").Value = '''5.550'
$ws.Range("E22").Value = '  +0.73%  '
$ws.Range("E23").Value = '  -0.37%  '
$ws.Range("D24").Value = '''6.556'
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("D25").Value = '''9.795'
$ws.Range("E25").Value = '  -0.66%  '
$ws.Range("D26").Value = '''166.31'
$ws.Range("E26").Value = '  -2.08%  '
$ws.Range("D27").Value = '''19.61'
$ws.Range("E27").Value = '  -3.49%  '
$ws.Range("D28").Value = '''2.160'
$ws.Range("E28").Value = '  -1.55%  '
$ws.Range("D29").Value = '''0.1038'
$ws.Range("E29").Value = '  -1.69%  '
$ws.Range("E30").Value = '  -3.23%  '
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '''4.580'
$ws.Range("E31").Value = '  -1.47%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''1.558'
$ws.Range("E32").Value = '  -2.27%  '
$ws.Range("D33").Value = '''4.388'
$ws.Range("E33").Value = '  -1.28%  '
$ws.Range("D34").Value = '''0.04884'
$ws.Range("E34").Value = '  -1.83%  '
$ws.Range("D35").Value = '''0.7565'
$ws.Range("E35").Value = '  -0.51%  '
$ws.Range("D36").Value = '''1.154'
$ws.Range("E36").Value = '  -1.31%  '
$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").Value = '''0.9996'
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").Value = '''2.740'
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").Value = '''0.01999'
$ws.Range("E39").Value = '  -2.08%  '
$ws.Range("D40").Value = '''79.03'
$ws.Range("E40").Value = '  +9.61%  '
$ws.Range("D41").Value = '''2.657'
$ws.Range("E41").Value = '  -2.14%  '
$ws.Range("D42").Value = '''6.459'
$ws.Range("E42").Value = '  +0.84%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''2.093'
$ws.Range("E43").Value = '  -5.06%  '
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").Value = '''0.9144'
$ws.Range("E44").Value = '  +3.65%  '
$ws.Range("D45").Value = '''0.4431'
$ws.Range("E45").Value = '  -2.73%  '
$ws.Range("D46").Value = '''107.91'
$ws.Range("E46").Value = '  -1.52%  '
$ws.Range("D47").Value = '''0.9990'
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("D48").Value = '''7.673'
$ws.Range("E48").Value = '  -2.75%  '
$ws.Range("D49").Value = '''983.77'
$ws.Range("E49").Value = '  +2.32%  '
$ws.Range("D50").Value = '''0.1244'
$ws.Range("E50").Value = '  -1.84%  '
$ws.Range("D51").Value = '''36.09'
$ws.Range("E51").Value = '  +0.88%  '
